$d = $word.ActiveDocument

# Helper: literal (non-wildcard) Find & ReplaceAll across the whole document
# story (main body). wdReplaceAll = 2, wdFindContinue wrap = 1
function Replace-All($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $replaceText, 2) | Out-Null
}

# --- Placeholder syntax migration: "<Placeholder>" -> "{{Placeholder}}" ---
# (a few also got their wording tidied up at the same time)

Replace-All "<Organisation Contact>" "{{Organisation Contact}}"
Replace-All "<Address from Organisation>" "{{Address from Organisation}}"
Replace-All "<External Ref>" "{{External Ref}}"
Replace-All "Our Ref: <CLO number>" "Our Ref: {{CLO number}}"
Replace-All "<LAG number>" "{{LAG number}}"
Replace-All "Contact: <Case Officer>" "Contact: {{Case Officer}}"
Replace-All "Direct Dial: <from case officer> " "Direct Dial: {{from case officer}} "
Replace-All "Email: <from Case Officer> " "Email: {{from Case Officer}} "
Replace-All "<Consultation Stage 1.10 Completion Date>" "{{Completion Date}}"
Replace-All "<FAO Organisation Contact>" "{{FAO Organisation Contact}}"
Replace-All "<Site Name>" "{{Site Name}}"
Replace-All "<Proposal>" "{{Proposal}}"
Replace-All "Thank you for your consultation received on <insert Log Date>." "Thank you for your consultation received on {{Log Date}}."
Replace-All "Local Plan: <insert name of selected APA>]" "Local Plan: {{insert name of selected APA}}]"
Replace-All "<Stage 1.10 Action = Predetermination: free text justification drawn from notes in stage 1.10>" "{{Stage 1.10 Action = Predetermination: free text justification drawn from notes in stage 1.10}}"
Replace-All "<Action for Consultation Stage 1.11 with standard scope note and any free text from notes>" "{{Action for Consultation Stage 1.11 with standard scope note and any free text from notes}}"
Replace-All "<Action for Consultation Stage 1.12 with standard scope note and any free text from notes>" "{{Action for Consultation Stage 1.12 with standard scope note and any free text from notes}}"
Replace-All "<Action for Consultation Stage 1.13 with standard scope note and any free text from notes>" "{{Action for Consultation Stage 1.13 with standard scope note and any free text from notes}}"
Replace-All "<Action for Consultation Stage 1.14 with standard scope note and any free text from notes>" "{{Action for Consultation Stage 1.14 with standard scope note and any free text from notes}}"
Replace-All "<insert digital signature>" "{{Digital Signature}}"
Replace-All "<Case Officer>" "{{Case Officer}}"

# --- Minor layout tweaks that came along with the same commit ---

# Page bottom margin nudged from 26.6pt (532 twips) to 26.65pt (533 twips)
$d.PageSetup.BottomMargin = 26.65

# The letterhead crest in the footer moves from "in front of text" to
# "behind text" (wp:anchor behindDoc 0 -> 1). Toggle the wrap type through
# wdWrapBehind and back to wdWrapTopAndBottom so only the behind/in-front
# flag flips and the original top/bottom wrap behaviour is preserved.
$sec = $d.Sections.First
$ftr = $sec.Footers.Item(1)
if ($ftr.Shapes.Count -gt 0) {
    $crest = $ftr.Shapes.Item(1)
    $crest.WrapFormat.Type = 5
    $crest.WrapFormat.Type = 4
}

Write-Output $d.Content.Text
